$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Apply updated coin data (price, link, name, volume) cell by cell to match the
# refreshed snapshot from the scheduled GitHub Actions run.

$ws.Range("D2").Value = '39.507.15'
$ws.Range("E2").Value = '  +1.67%  '

$ws.Range("D3").Value = '2.157.27'
$ws.Range("E3").Value = '  +3.08%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '229.23'
$ws.Range("E5").Value = '  +0.24%  '

$ws.Range("E6").Value = '  +0.94%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '63.02'
$ws.Range("E7").Value = '  +4.08%  '

$ws.Range("E8").Value = '  -0.03%  '

$ws.Range("E9").Value = '  +2.66%  '

$ws.Range("E10").Value = '  +3.00%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.104'
$ws.Range("E11").Value = '  -0.15%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '16.15'
$ws.Range("E12").Value = '  +7.82%  '

$ws.Range("D13").Value = '2.475.67'
$ws.Range("E13").Value = '  +3.02%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '22.25'
$ws.Range("E14").Value = '  +1.56%  '

$ws.Range("E15").Value = '  +3.44%  '

$ws.Range("E16").Value = '  +1.90%  '

$ws.Range("D17").Value = '2.228.92'
$ws.Range("E17").Value = '  +6.29%  '

$ws.Range("D18").Value = '39.474.45'
$ws.Range("E18").Value = '  +1.96%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '72.36'
$ws.Range("E19").Value = '  +1.05%  '

$ws.Range("E20").Value = '  +1.80%  '

$ws.Range("E21").Value = '  +2.11%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '228.52'
$ws.Range("E22").Value = '  +0.56%  '

$ws.Range("E23").Value = '  +0.00%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.35'
$ws.Range("E24").Value = '  -1.31%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.37'
$ws.Range("E25").Value = '  +1.00%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.77'
$ws.Range("E26").Value = '  +3.01%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '172.43'
$ws.Range("E27").Value = '  +0.66%  '

$ws.Range("E28").Value = '  -0.66%  '

$ws.Range("E29").Value = '  -3.93%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.65'
$ws.Range("E30").Value = '  +2.50%  '

$ws.Range("E31").Value = '  +8.23%  '

$ws.Range("E32").Value = '  +1.34%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.64'
$ws.Range("E33").Value = '  +3.02%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.81'
$ws.Range("E34").Value = '  +2.41%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '7.13'
$ws.Range("E35").Value = '  +10.64%  '

$ws.Range("E36").Value = '  +1.93%  '

$ws.Range("E37").Value = '  +2.61%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.56'
$ws.Range("E38").Value = '  -0.30%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.999'
$ws.Range("E39").Value = '  -0.20%  '

$ws.Range("E40").Value = '  +3.34%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '18.18'
$ws.Range("E41").Value = '  +0.17%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '103.37'
$ws.Range("E42").Value = '  +2.42%  '

$ws.Range("D43").Value = '1.531.87'
$ws.Range("E43").Value = '  -0.61%  '

$ws.Range("E44").Value = '  +5.43%  '

$ws.Range("B45").Value = 'FTXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.36'
$ws.Range("E45").Value = '  +5.95%  '

$ws.Range("B46").Value = 'Cronos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0927'
$ws.Range("E46").Value = '  +0.50%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.11'
$ws.Range("E47").Value = '  +6.83%  '

$ws.Range("E48").Value = '  -0.32%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.75'
$ws.Range("E49").Value = '  +1.39%  '

$ws.Range("D50").Value = '2.359.54'
$ws.Range("E50").Value = '  +3.08%  '

$ws.Range("E51").Value = '  +0.05%  '
